# Test Data Generator - Deirdre.xlsx edit script
# Implements: clear the two hard-coded "insert into user_profile" demo
# formulas on user_profile, and add three new reference sheets
# (relationship_type, user_interests, interests) with their header/data rows.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- user_profile: drop the two illustrative formula cells in L1/L2 ---
# L1 keeps its bold/wrap style (s=5) but becomes blank.
$ws1.Range("L1").ClearContents() | Out-Null
# L2 is fully cleared (style + content), matching the removed <c> element.
$ws1.Range("L2").Clear() | Out-Null

# Row heights shrink now that the wrapped formula text is gone.
$ws1.Rows.Item(1).RowHeight = 19.5
$ws1.Rows.Item(2).AutoFit() | Out-Null

# --- add the three new worksheets, in tab order after user_profile ---
$wsRel = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsRel.Name = "relationship_type"

$wsUI = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsUI.Name = "user_interests"

$wsInt = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsInt.Name = "interests"

# --- relationship_type sheet ---
$wsRel.Range("A1").Value = "User_id"
$wsRel.Range("B1").Value = "relationship_type"
$wsRel.Range("C1").Value = " "
$wsRel.Range("D1").Value = " "
$wsRel.Range("E1").Value = " "
$wsRel.Range("B2").Value = "love"
$wsRel.Range("B3").Value = "casual"
$wsRel.Range("B4").Value = "friendship"
$wsRel.Range("B5").Value = "relationship"

$ws1.Range("A1").Copy() | Out-Null
$wsRel.Range("A1:E1").PasteSpecial(-4122) | Out-Null
$wsRel.Application.CutCopyMode = $false

$wsRel.Columns.Item(2).ColumnWidth = 36.5
$wsRel.Rows.Item(1).RowHeight = 19.5
$wsRel.Range("B1").Select() | Out-Null

# --- user_interests sheet ---
$wsUI.Range("B1").Value = "interest_id"
$wsUI.Range("A1").Value = "user_id"
$wsUI.Range("B2").Value = 1
$wsUI.Range("B3").Value = 2
$wsUI.Range("B4").Value = 3
$wsUI.Range("B5").Value = 4
$wsUI.Range("B6").Value = 5

$ws1.Range("A1").Copy() | Out-Null
$wsUI.Range("A1:B1").PasteSpecial(-4122) | Out-Null
$wsUI.Application.CutCopyMode = $false

$wsUI.Columns.Item(2).ColumnWidth = 19.333333333333332
$wsUI.Rows.Item(1).RowHeight = 19.5
$wsUI.Rows.Item(1).Select() | Out-Null

# --- interests sheet ---
$wsInt.Range("A1").Value = "user_id"
$wsInt.Range("B1").Value = "interest_id"
$wsInt.Range("C1").Value = "description"

$ws1.Range("A1").Copy() | Out-Null
$wsInt.Range("A1:C1").PasteSpecial(-4122) | Out-Null
$wsInt.Application.CutCopyMode = $false

$wsInt.Columns.Item(2).ColumnWidth = 14.333333333333316
$wsInt.Columns.Item(3).ColumnWidth = 15.333333333333313
$wsInt.Rows.Item(1).RowHeight = 19.5
$wsInt.Range("C1").Select() | Out-Null

# Set sheet1's lingering selection (it was last viewed scrolled to L7)
# before leaving it, then land on user_interests - the sheet active when
# the workbook was saved.
$ws1.Activate()
$ws1.Range("L7").Select() | Out-Null
$wsUI.Activate()
